# "Change shape of the notes" - mark the requirement about customizable
# note shapes as struck-through (descoped), matching the styling already
# used on the other removed/descoped bullet points in this list.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "forma geometrica a notelor din piano roll sa poata fi customizata") {
        $p.Range.Font.StrikeThrough = $true
    }
}
